$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.799.25"
$ws.Range("E2").Value = "'  +2.46%  "

$ws.Range("D3").Value = "'3.951.68"
$ws.Range("E3").Value = "'  +1.26%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'528.14"
$ws.Range("E5").Value = "'  +8.43%  "

$ws.Range("D6").Value = "'146.87"
$ws.Range("E6").Value = "'  +0.57%  "

$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "'  +0.32%  "

$ws.Range("E8").Value = "'  -0.08%  "

$ws.Range("D9").Value = "'0.729"
$ws.Range("E9").Value = "'  +0.48%  "

$ws.Range("D10").Value = "'0.174"
$ws.Range("E10").Value = "'  +5.83%  "

$ws.Range("D11").Value = "'0.0000343"
$ws.Range("E11").Value = "'  +0.49%  "

$ws.Range("D12").Value = "'42.89"
$ws.Range("E12").Value = "'  +0.39%  "

$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.46"
$ws.Range("E13").Value = "'  -2.09%  "

$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'4.581.40"
$ws.Range("E14").Value = "'  +1.16%  "

$ws.Range("D15").Value = "'3.957.51"
$ws.Range("E15").Value = "'  +1.05%  "

$ws.Range("D16").Value = "'14.22"
$ws.Range("E16").Value = "'  +1.45%  "

$ws.Range("E17").Value = "'  -0.03%  "

$ws.Range("E18").Value = "'  +7.72%  "

$ws.Range("D19").Value = "'19.87"
$ws.Range("E19").Value = "'  +0.53%  "

$ws.Range("D20").Value = "'69.782.41"
$ws.Range("E20").Value = "'  +2.15%  "

$ws.Range("D21").Value = "'436.75"
$ws.Range("E21").Value = "'  +0.34%  "

$ws.Range("D22").Value = "'3.40"
$ws.Range("E22").Value = "'  -3.38%  "

$ws.Range("D23").Value = "'14.54"
$ws.Range("E23").Value = "'  -2.16%  "

$ws.Range("D24").Value = "'88.55"
$ws.Range("E24").Value = "'  +1.28%  "

$ws.Range("D25").Value = "'4.05"
$ws.Range("E25").Value = "'  +11.96%  "

$ws.Range("D26").Value = "'11.99"
$ws.Range("E26").Value = "'  +7.53%  "

$ws.Range("D27").Value = "'11.08"
$ws.Range("E27").Value = "'  -1.81%  "

$ws.Range("D28").Value = "'36.82"
$ws.Range("E28").Value = "'  -3.18%  "

$ws.Range("D30").Value = "'699.54"
$ws.Range("E30").Value = "'  -3.49%  "

$ws.Range("D31").Value = "'13.38"
$ws.Range("E31").Value = "'  -1.78%  "

$ws.Range("E32").Value = "'  -1.93%  "

$ws.Range("E33").Value = "'  -0.91%  "

$ws.Range("D34").Value = "'67.12"
$ws.Range("E34").Value = "'  +11.78%  "

$ws.Range("D35").Value = "'0.444"
$ws.Range("E35").Value = "'  +9.97%  "

$ws.Range("D36").Value = "'0.0" + [char]0x2083 + "0877"
$ws.Range("E36").Value = "'  +2.91%  "

$ws.Range("D37").Value = "'5.98"
$ws.Range("E37").Value = "'  -3.68%  "

$ws.Range("D38").Value = "'40.38"
$ws.Range("E38").Value = "'  -2.55%  "

$ws.Range("D39").Value = "'0.149"
$ws.Range("E39").Value = "'  +1.39%  "

$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "'  +0.17%  "

$ws.Range("E41").Value = "'  -0.10%  "

$ws.Range("E42").Value = "'  +1.40%  "

$ws.Range("B43").Value = "'Fetch.AI"
$ws.Range("C43").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.86"
$ws.Range("E43").Value = "'  -2.66%  "

$ws.Range("B44").Value = "'WEMIXToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'3.14"
$ws.Range("E44").Value = "'  +7.55%  "

$ws.Range("E45").Value = "'  -3.86%  "

$ws.Range("E46").Value = "'  +1.26%  "

$ws.Range("D47").Value = "'3.39"
$ws.Range("E47").Value = "'  +3.56%  "

$ws.Range("D48").Value = "'0.0" + [char]0x2086 + "0364"
$ws.Range("E48").Value = "'  +11.01%  "

$ws.Range("D49").Value = "'3.07"
$ws.Range("E49").Value = "'  +9.55%  "

$ws.Range("D50").Value = "'3.36"
$ws.Range("E50").Value = "'  -0.70%  "

$ws.Range("D51").Value = "'2.10"
$ws.Range("E51").Value = "'  -1.00%  "
